$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of data for sg_rr_68_025 2023-12-11 17-15-27 FWHM run
$ws.Cells.Item(57, 1).Value = "sg_rr_68_025 2023-12-11 17-15-27.csv"
$ws.Cells.Item(57, 2).Value = 0.01
$ws.Cells.Item(57, 3).Value = 1000
$ws.Cells.Item(57, 4).Value = 5001
$ws.Cells.Item(57, 5).Value = 1530
$ws.Cells.Item(57, 6).Value = 1570
$ws.Cells.Item(57, 7).Value = 0.5
$ws.Cells.Item(57, 8).Value = "(approx_fsr/2)/wavelength step size"
$ws.Cells.Item(57, 9).Value = 1.7
$ws.Cells.Item(57, 10).Value = 1.4480769230769099
$ws.Cells.Item(57, 11).Value = 0.00622696940163916
$ws.Cells.Item(57, 12).Value = "yes (although start peak maybe missed)"
$ws.Cells.Item(57, 13).Value = 0.137531855376269
$ws.Cells.Item(57, 14).Value = 0.0020576576158408901

# Adjust column L width to match the diff (Excel stores width in a font-derived
# pixel unit; 7.1 "characters" round-trips to the target stored width of 8)
$ws.Columns.Item(12).ColumnWidth = 7.1

# Scroll the view and reselect so the view state matches the new row count
$excel.ActiveWindow.ScrollRow = 43
$ws.Range("A57").Select() | Out-Null
